# "Slight changes to calculator"
# - TGDraw: Relays device count 1 -> 0
# - TGDraw: Camera device count 1 -> 2
# - TGDraw: J12 formula doubled (Sources!F3/TGDraw!J11 -> 2*Sources!F3/TGDraw!J11)
# - View state: TGDraw becomes the active/selected sheet (was TBDraw),
#   with new selections on both sheets.

$wb = $excel.ActiveWorkbook

$tb = $wb.Worksheets.Item("TBDraw")
$tg = $wb.Worksheets.Item("TGDraw")

# Device count changes on TGDraw
$tg.Range("F5").Value = 0
$tg.Range("F7").Value = 2

# Double the TG battery run-time formula
$tg.Range("J12").Formula = "=2*Sources!F3/TGDraw!J11"

# Update selections, then make TGDraw the active sheet/tab
$tb.Range("A24").Select()

$tg.Activate()
$tg.Range("J13").Select()

$wb.Save()
